$d = $word.ActiveDocument

function New-WordXmlFragment([string]$InnerBodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $InnerBodyXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# 1) Remove the existing "_GoBack" bookmark from its old spot; it will be
#    recreated at the end of the new "Python 3.8.0" paragraph below.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2) "<구조>" (paragraph 1) becomes "<사용 버전>"
$d.Content.Find.Execute("<구조>", $true, $false, $false, $false, $false, $true, 1, $false, "<사용 버전>", 2) | Out-Null

# 3) Insert the new "Python 3.8.0" paragraph (with the relocated bookmark)
#    right after the (now renamed) first paragraph.
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertParagraphAfter() | Out-Null
$pythonPara = $d.Paragraphs.Item(2)
$pythonXml = New-WordXmlFragment('<w:body><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>P</w:t></w:r><w:r><w:t>ython 3.8.0</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body>')
$pythonPara.Range.InsertXML($pythonXml) | Out-Null

# 4) Insert a blank paragraph after it.
$pythonPara2 = $d.Paragraphs.Item(2)
$pythonPara2.Range.InsertParagraphAfter() | Out-Null
$blankPara = $d.Paragraphs.Item(3)
$blankXml = New-WordXmlFragment('<w:body><w:p/></w:body>')
$blankPara.Range.InsertXML($blankXml) | Out-Null

# 5) Insert the relocated "<구조>" heading paragraph after the blank one.
$blankPara2 = $d.Paragraphs.Item(3)
$blankPara2.Range.InsertParagraphAfter() | Out-Null
$structPara = $d.Paragraphs.Item(4)
$structXml = New-WordXmlFragment('<w:body><w:p><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>&lt;구조&gt;</w:t></w:r></w:p></w:body>')
$structPara.Range.InsertXML($structXml) | Out-Null
